# Sample Project / Main.xlsx — project saved.
# Rule row 11 ("R40" rule row): the Rule-name cell B11 is changed from the
# text "R40" to the text "1" (still a text value, stored as a shared
# string - not a number). A leading apostrophe forces Excel to keep the
# numeric-looking entry as literal text instead of coercing it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("B11").Value = "'1"
